$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(3, 4).Range.Text = "0.05"
$t.Cell(4, 4).Range.Text = "0.24"
$t.Cell(5, 4).Range.Text = "0.10"
$t.Cell(6, 4).Range.Text = "0.61"
$t.Cell(7, 4).Range.Text = "0.03"
$t.Cell(8, 4).Range.Text = "0.18"
$t.Cell(9, 4).Range.Text = "0.10"
$t.Cell(10, 4).Range.Text = "0.69"
$t.Cell(11, 4).Range.Text = "0.14"
$t.Cell(14, 4).Range.Text = "0.44"
$t.Cell(15, 4).Range.Text = "0.19"
$t.Cell(19, 4).Range.Text = "0.43"
$t.Cell(20, 4).Range.Text = "0.12"
$t.Cell(21, 4).Range.Text = "0.11"
$t.Cell(22, 4).Range.Text = "0.33"
$t.Cell(23, 4).Range.Text = "0.48"
$t.Cell(24, 4).Range.Text = "0.07"
$t.Cell(25, 4).Range.Text = "0.05"
$t.Cell(26, 4).Range.Text = "0.39"
$t.Cell(28, 4).Range.Text = "0.02"
$t.Cell(29, 4).Range.Text = "0.61"
$t.Cell(30, 4).Range.Text = "0.14"
$t.Cell(31, 4).Range.Text = "0.69"
$t.Cell(33, 4).Range.Text = "0.20"
$t.Cell(35, 4).Range.Text = "0.07"
$t.Cell(36, 4).Range.Text = "0.43"
$t.Cell(37, 4).Range.Text = "0.24"
$t.Cell(38, 4).Range.Text = "0.26"
$t.Cell(39, 4).Range.Text = "0.04"
$t.Cell(40, 4).Range.Text = "0.53"
$t.Cell(41, 4).Range.Text = "0.38"
$t.Cell(42, 4).Range.Text = "0.06"
$t.Cell(43, 4).Range.Text = "0.05"
$t.Cell(44, 4).Range.Text = "0.35"
$t.Cell(45, 4).Range.Text = "0.39"
$t.Cell(46, 4).Range.Text = "0.21"
$t.Cell(47, 4).Range.Text = "0.16"
$t.Cell(48, 4).Range.Text = "0.26"
$t.Cell(49, 4).Range.Text = "0.25"
$t.Cell(50, 4).Range.Text = "0.34"
$t.Cell(51, 4).Range.Text = "0.18"
$t.Cell(52, 4).Range.Text = "0.23"
$t.Cell(53, 4).Range.Text = "0.35"
$t.Cell(54, 4).Range.Text = "0.24"
$t.Cell(55, 4).Range.Text = "0.15"
$t.Cell(56, 4).Range.Text = "0.20"
$t.Cell(57, 4).Range.Text = "0.34"
$t.Cell(58, 4).Range.Text = "0.31"
$t.Cell(59, 4).Range.Text = "0.07"
$t.Cell(60, 4).Range.Text = "0.24"
$t.Cell(61, 4).Range.Text = "0.48"
$t.Cell(62, 4).Range.Text = "0.21"
$t.Cell(63, 4).Range.Text = "0.07"
$t.Cell(64, 4).Range.Text = "0.46"
$t.Cell(65, 4).Range.Text = "0.20"
$t.Cell(66, 4).Range.Text = "0.27"
$t.Cell(67, 4).Range.Text = "0.03"
$t.Cell(68, 4).Range.Text = "0.06"
$t.Cell(69, 4).Range.Text = "0.55"
$t.Cell(70, 4).Range.Text = "0.35"
$t.Cell(71, 4).Range.Text = "0.02"
$t.Cell(72, 4).Range.Text = "0.06"
$t.Cell(73, 4).Range.Text = "0.56"
$t.Cell(74, 4).Range.Text = "0.36"
